$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# 1) Update the "Förändrad" date column (C) from 45207 to 45208 for every data row (2-28)
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}

# 2) Update every HYPERLINK formula so its folder segment changes
#    from "Logging_HEBY" to "Logging_0331" (columns S through Y)
$cols = @(19, 20, 21, 22, 23, 24, 25)  # S, T, U, V, W, X, Y
for ($r = 2; $r -le 28; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -and $f -like "*Logging_HEBY*") {
            $cell.Formula = $f.Replace("Logging_HEBY", "Logging_0331")
        }
    }
}
